# Add new GlobalConstant rows (Ev13CountLimit .. Ev17CountLimit) to the
# "GlobalConstantIntTable" worksheet, which is the workbook's active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newConstants = @(
    @{ Row = 47; Name = "Ev13CountLimit"; Value = 4 },
    @{ Row = 48; Name = "Ev14CountLimit"; Value = 8 },
    @{ Row = 49; Name = "Ev15CountLimit"; Value = 9 },
    @{ Row = 50; Name = "Ev16CountLimit"; Value = 12 },
    @{ Row = 51; Name = "Ev17CountLimit"; Value = 17 }
)

foreach ($item in $newConstants) {
    $ws.Cells.Item($item.Row, 1).Value = $item.Name
    $ws.Cells.Item($item.Row, 2).Value = $item.Value
}
